# Add two new columns (I: "I0", J: "IF") to the worksheet with header
# labels in row 1 and numeric data in rows 2-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - match style of existing headers (B1:H1 use style index 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Numeric data for rows 2-31 (columns I and J)
$data = @{
    2  = @(11, 11)
    3  = @(6, 6)
    4  = @(7, 7)
    5  = @(7, 7)
    6  = @(7, 7)
    7  = @(8, 9)
    8  = @(8, 8)
    9  = @(4, 5)
    10 = @(9, 10)
    11 = @(9, 9)
    12 = @(7, 7)
    13 = @(5, 5)
    14 = @(7, 8)
    15 = @(3, 5)
    16 = @(3, 3)
    17 = @(6, 7)
    18 = @(9, 9)
    19 = @(7, 7)
    20 = @(4, 5)
    21 = @(6, 7)
    22 = @(7, 7)
    23 = @(3, 3)
    24 = @(7, 7)
    25 = @(3, 4)
    26 = @(9, 9)
    27 = @(7, 7)
    28 = @(4, 4)
    29 = @(4, 4)
    30 = @(5, 5)
    31 = @(4, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
